# Replace the placeholder homework answer with the real one, and move the
# "_GoBack" bookmark (which Word drops at the site of the last edit) from
# the trailing empty paragraph to right after the newly-typed text.

$d = $word.ActiveDocument

# 1. Swap the TODO placeholder text for the actual answer.
$d.Content.Find.Execute("TODO: FINISH", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "a*b*c*, a*b*c, a*bc+", 2)

# 2. Drop the old "_GoBack" bookmark that currently sits in the last
#    (empty) paragraph.
$d.Bookmarks.Item("_GoBack").Delete()

# 3. Re-create "_GoBack" as a zero-length bookmark right after the text
#    we just inserted (i.e. before the line break that ends that
#    paragraph), matching where Word leaves it after a real edit.
$r = $d.Content
$r.Find.Execute("a*b*c*, a*b*c, a*bc+", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
